$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

$ws.Range("M8").Value = "سيرفيس "
$ws.Range("N8").Value = "تم تغير اول جريده وتغير فلاتس متحركه"

$ws.Range("A31").Value = "nan"
$ws.Range("L31").Value = "nan"
$ws.Range("M31").Value = "nan"
$ws.Range("N31").Value = "nan"
